$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159; this shifts existing rows 159:260 down to
# 160:261, carrying over their values/styles (matches the diff, where every
# row N (160..261) now holds what used to be row N-1, and a brand-new record
# lands at the top in row 159).
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with this week's new record.
$ws.Cells.Item(159, 1).Value = 11
$ws.Cells.Item(159, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(159, 3).Value = "Bíobío"
$ws.Cells.Item(159, 4).Value = 45001
$ws.Cells.Item(159, 5).Value = 8
$ws.Cells.Item(159, 6).Value = 100112003
$ws.Cells.Item(159, 7).Value = "Ajo"
$ws.Cells.Item(159, 8).Value = "Chino"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 220
$ws.Cells.Item(159, 11).Value = 15000
$ws.Cells.Item(159, 12).Value = 16000
$ws.Cells.Item(159, 13).Value = 15455
$ws.Cells.Item(159, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(159, 15).Value = "China"
$ws.Cells.Item(159, 16).Value = 1546
$ws.Cells.Item(159, 17).Value = 10
$ws.Cells.Item(159, 18).Value = "Hortaliza"
